# Actualización automática del mapa (2025-08-07 14:26:25)
# Adds the new incident row (row 61) to the "NEW" worksheet, mirroring the
# data appended by the automated map-update job.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRow = $ws.UsedRange.Rows.Count + 1

# Columns A-E contain values that look like numbers or dates ("6579",
# "8/7/2025", "13", "808749184") but must be stored as literal text, exactly
# like the other rows in this sheet. Pre-formatting the cells as Text before
# assigning the values prevents Excel from auto-converting them into numbers
# or date serials.
$textRange = $ws.Range("A" + $newRow + ":E" + $newRow)
$textRange.NumberFormat = "@"

$ws.Range("A" + $newRow).Value = "6579"
$ws.Range("B" + $newRow).Value = "8/7/2025"
$ws.Range("C" + $newRow).Value = "RIVADAVIA MARTIN, COMODORO 1350"
$ws.Range("D" + $newRow).Value = "13"
$ws.Range("E" + $newRow).Value = "808749184"
$ws.Range("F" + $newRow).Value = "NEW"
$ws.Range("G" + $newRow).Value = "Pendiente"
$ws.Range("H" + $newRow).Value = "Poste inclinado"
$ws.Range("I" + $newRow).Value = 1
$ws.Range("J" + $newRow).Value = "Aplomo"
$ws.Range("K" + $newRow).Value = "Sin equipos"
$ws.Range("L" + $newRow).Value = "Poste"
$ws.Range("M" + $newRow).Value = -58.461024
$ws.Range("N" + $newRow).Value = -34.539409
$ws.Range("O" + $newRow).Value = "Saavedra"
$ws.Range("P" + $newRow).Value = "Capital Norte"
